$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 14 (keep header row 1 and data rows 2-5)
$ws.Range("A6:B14").EntireRow.Delete() | Out-Null

# Update the remaining data rows (2-5) with the new bin labels and values
$ws.Range("A2").Value = "10-15"
$ws.Range("B2").Value = 3.411006674596242

$ws.Range("A3").Value = "5-10"
$ws.Range("B3").Value = 3.003127636795952

$ws.Range("A4").Value = "10-15"
$ws.Range("B4").Value = 3.925750268684639

$ws.Range("A5").Value = "5-10"
$ws.Range("B5").Value = 3.350359798132719
